# Updated maxmsgsize for grpc
# Re-populate the unit test data rows (A2:F25) on the active sheet to match
# the refreshed test fixture values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $values) {
    $arr = New-Object 'object[,]' 1, $values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range("A$($row):F$($row)").Value = $arr
}

Set-RowValues 2  @(201,  9, 30, 15, 45, 30)
Set-RowValues 3  @(801,  3, 67, 65, 52, 45)
Set-RowValues 4  @(901, 16, 15, 45, 60, 60)
Set-RowValues 5  @(301,  6, 45, 30, 60, 45)
Set-RowValues 6  @(1202, 2, 10, 10, 10, 10)
Set-RowValues 7  @(902,  1,  0,  0,  0,  0)
Set-RowValues 8  @(1001,18, 30, 75, 60, 72)
Set-RowValues 9  @(401,  9, 48, 67, 75, 45)
Set-RowValues 10 @(701,  3, 90, 45, 97, 15)
Set-RowValues 11 @(1201, 2, 10, 10, 10, 10)
Set-RowValues 12 @(101,  9, 30, 15, 60, 15)
Set-RowValues 13 @(501,  9, 52, 30, 75, 45)
Set-RowValues 14 @(601,  9, 60, 67, 60, 42)
Set-RowValues 16 @(3,    0,  3,  3,  3,  3)
Set-RowValues 17 @(2,    0,  2,  2,  2,  2)
Set-RowValues 19 @(802,  0,  4,  5,  4,  0)
Set-RowValues 20 @(1101, 0, 15, 30, 30,  0)
Set-RowValues 21 @(1,    0,  2,  2,  2,  2)
Set-RowValues 22 @(402,  0,  0,  4,  0,  0)
Set-RowValues 23 @(602,  0,  0,  4,  0,  9)
